$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Change 1: highlight (darkRed) the paragraph
# "A bíblia a mensagem ela interpreta os textos bíblicos e não traduz a
# mensagem." - both the paragraph mark (pPr/rPr) and the run text need
# the highlight applied.
# ----------------------------------------------------------------------
$searchText1 = "A bíblia a mensagem ela interpreta os textos bíblicos e não traduz a mensagem."
$probe1 = $d.Content.Duplicate
$found1 = $probe1.Find.Execute($searchText1, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $targetStart1 = $probe1.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -eq $targetStart1) {
            $p.Range.Font.HighlightColorIndex = 13   # wdDarkRed
            break
        }
    }
}

# ----------------------------------------------------------------------
# Change 2: fix the typo "numero" -> "número" (only the first occurrence,
# right after "O "), and drop the gramStart/gramEnd proofErr markers
# that used to bracket it. The whole host paragraph is rebuilt (keeping
# its original run split / formatting) via InsertXML so that the
# now-orphaned proofErr markers are dropped cleanly.
# ----------------------------------------------------------------------
$probe2 = $d.Content.Duplicate
$found2 = $probe2.Find.Execute("numero", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $wordStart = $probe2.Start

    $hostParaStart = -1
    $hostParaEnd = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($wordStart -ge $p.Range.Start -and $wordStart -lt $p.Range.End) {
            $hostParaStart = $p.Range.Start
            $hostParaEnd = $p.Range.End
            break
        }
    }

    $target = $d.Range($hostParaStart, $hostParaEnd)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">O </w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>número</w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> 12 na numerologia Bíblica representa a Teocracia, representa o numero da escolha de Deus, AP 21 eram 12 colunas, 12 portas, 24 anciãos (Representa os apóstolos e as 12 tribos e os apóstolos a igreja adorando ao cordeiro), 144 tribos.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)
}

Write-Host "Done"
